$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.363.67'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').Value = '1.909.99'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.730'
$ws.Range('E5').Value = '  +10.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '256.34'
$ws.Range('E6').Value = '  +4.31%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.91'
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.370'
$ws.Range('E9').Value = '  +6.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '53.17'
$ws.Range('E10').Value = '  +0.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0760'
$ws.Range('E11').Value = '  +6.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0989'
$ws.Range('E12').Value = '  -0.24%  '
$ws.Range('D13').Value = '2.186.15'
$ws.Range('E13').Value = '  +0.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.94'
$ws.Range('E14').Value = '  +6.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.734'
$ws.Range('E15').Value = '  +5.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.99'
$ws.Range('E16').Value = '  +4.30%  '
$ws.Range('D17').Value = '1.904.97'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('D18').Value = '35.333.79'
$ws.Range('E18').Value = '  +0.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '75.32'
$ws.Range('E19').Value = '  +4.33%  '
$ws.Range('D20').Value = '0.0₃0847'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '245.20'
$ws.Range('E21').Value = '  +2.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.12'
$ws.Range('E22').Value = '  +5.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.14'
$ws.Range('E23').Value = '  +7.08%  '
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.45'
$ws.Range('E25').Value = '  +7.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.46'
$ws.Range('E26').Value = '  +3.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.28'
$ws.Range('E27').Value = '  -2.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.76'
$ws.Range('E28').Value = '  +3.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.84'
$ws.Range('E29').Value = '  +3.07%  '
$ws.Range('E30').Value = '  +5.29%  '
$ws.Range('D31').Value = '4.127.35'
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('E32').Value = '  +6.49%  '
$ws.Range('E33').Value = '  +24.29%  '
$ws.Range('E34').Value = '  +14.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0591'
$ws.Range('E35').Value = '  +5.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.28'
$ws.Range('E36').Value = '  +5.28%  '
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.918'
$ws.Range('E38').Value = '  -2.11%  '
$ws.Range('E39').Value = '  +1.34%  '
$ws.Range('E40').Value = '  +6.83%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '97.45'
$ws.Range('E41').Value = '  +9.03%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.08'
$ws.Range('E42').Value = '  +6.52%  '
$ws.Range('E43').Value = '  +3.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0645'
$ws.Range('E44').Value = '  +2.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.48'
$ws.Range('E45').Value = '  +4.77%  '
$ws.Range('D46').Value = '1.338.64'
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.74'
$ws.Range('E48').Value = '  +4.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.76'
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.14'
$ws.Range('E50').Value = '  -7.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0757'
$ws.Range('E51').Value = '  +6.84%  '
